$wb = $excel.ActiveWorkbook

# The "comps" tab is currently the last sheet and is the selected tab.
# Add a new worksheet right after it, named "spaceInColHeader".
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "spaceInColHeader"

# Populate the new sheet: header row with a trailing-space column name,
# plus two data rows.
$ws.Range("A1").Value = "col1"
$ws.Range("B1").Value = "col2 "
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 4

# Match the selection state recorded for the new sheet.
$ws.Range("B1").Select() | Out-Null
